$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, $Addr, $Val)
    $c = $Sheet.Range($Addr)
    $c.NumberFormat = "@"
    $c.Value = $Val
    $c.NumberFormat = "General"
}

Set-TextValue $ws 'D2' '246.75'
Set-TextValue $ws 'G2' '20'
Set-TextValue $ws 'D3' '22.39'
Set-TextValue $ws 'G3' '20'
Set-TextValue $ws 'D4' '5.242'
Set-TextValue $ws 'G4' '20'
Set-TextValue $ws 'D5' '0.05681'
Set-TextValue $ws 'G5' '20'
Set-TextValue $ws 'G6' '20'
Set-TextValue $ws 'D7' '6.310'
Set-TextValue $ws 'G7' '20'
Set-TextValue $ws 'D8' '0.8071'
Set-TextValue $ws 'G8' '20'
Set-TextValue $ws 'D9' '0.8731'
Set-TextValue $ws 'G9' '20'
Set-TextValue $ws 'D10' '0.1413'
Set-TextValue $ws 'G10' '20'
Set-TextValue $ws 'D11' '0.07403'
Set-TextValue $ws 'G11' '20'
Set-TextValue $ws 'D12' '0.03035'
Set-TextValue $ws 'G12' '20'
Set-TextValue $ws 'D13' '0.03074'
Set-TextValue $ws 'G13' '20'
Set-TextValue $ws 'D14' '0.09400'
Set-TextValue $ws 'G14' '20'
Set-TextValue $ws 'D15' '3.886'
Set-TextValue $ws 'G15' '20'
Set-TextValue $ws 'D16' '0.001574'
Set-TextValue $ws 'G16' '20'
Set-TextValue $ws 'D17' '0.04784'
Set-TextValue $ws 'G17' '20'
Set-TextValue $ws 'B18' 'TigerCash'
Set-TextValue $ws 'C18' 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws 'D18' '0.006382'
Set-TextValue $ws 'E18' '17TigerCashTCH'
Set-TextValue $ws 'G18' '20'
Set-TextValue $ws 'B19' 'HotbitToken'
Set-TextValue $ws 'C19' 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
Set-TextValue $ws 'D19' '0.005032'
Set-TextValue $ws 'E19' '18HotbitTokenHTB'
Set-TextValue $ws 'G19' '20'
Set-TextValue $ws 'B20' 'BitKan'
Set-TextValue $ws 'C20' 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
Set-TextValue $ws 'D20' '0.0009967'
Set-TextValue $ws 'E20' '19BitKanKAN'
Set-TextValue $ws 'G20' '20'
Set-TextValue $ws 'B21' 'NitroEx'
Set-TextValue $ws 'C21' 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
Set-TextValue $ws 'D21' '0.0001501'
Set-TextValue $ws 'E21' '20NitroExNTX'
Set-TextValue $ws 'G21' '20'
Set-TextValue $ws 'B22' 'LEO'
Set-TextValue $ws 'C22' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws 'D22' '3.689'
Set-TextValue $ws 'E22' '21LEOLEO'
Set-TextValue $ws 'G22' '20'
Set-TextValue $ws 'B23' 'BTSEToken'
Set-TextValue $ws 'C23' 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue $ws 'D23' '2.199'
Set-TextValue $ws 'E23' '22BTSETokenBTSE'
Set-TextValue $ws 'G23' '20'
Set-TextValue $ws 'B24' 'One'
Set-TextValue $ws 'C24' 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue $ws 'D24' '0.01090'
Set-TextValue $ws 'E24' '23OneONEBestin24h'
Set-TextValue $ws 'G24' '20'
Set-TextValue $ws 'D25' '0.3279'
Set-TextValue $ws 'G25' '20'
Set-TextValue $ws 'G26' '20'
Set-TextValue $ws 'D27' '0.0004753'
Set-TextValue $ws 'E27' '26UpBotsUBXT'
Set-TextValue $ws 'G27' '20'
Set-TextValue $ws 'G28' '20'
Set-TextValue $ws 'G29' '20'
Set-TextValue $ws 'G30' '20'
Set-TextValue $ws 'G31' '20'
Set-TextValue $ws 'G32' '20'
Set-TextValue $ws 'G33' '20'
Set-TextValue $ws 'G34' '20'
Set-TextValue $ws 'G35' '20'
Set-TextValue $ws 'G36' '20'
Set-TextValue $ws 'G37' '20'
Set-TextValue $ws 'G38' '20'
Set-TextValue $ws 'G39' '20'
Set-TextValue $ws 'D40' '0.03943'
Set-TextValue $ws 'G40' '20'
Set-TextValue $ws 'B41' 'KickToken'
Set-TextValue $ws 'C41' 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue $ws 'D41' '0.006804'
Set-TextValue $ws 'E41' '40KickTokenKICK'
Set-TextValue $ws 'G41' '20'
Set-TextValue $ws 'B42' 'BKEXToken'
Set-TextValue $ws 'C42' 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue $ws 'D42' '0.1065'
Set-TextValue $ws 'E42' '41BKEXTokenBKK'
Set-TextValue $ws 'G42' '20'
Set-TextValue $ws 'B43' 'CEJI'
Set-TextValue $ws 'C43' 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue $ws 'D43' '0.003202'
Set-TextValue $ws 'E43' '42CEJICEJI'
Set-TextValue $ws 'G43' '20'
Set-TextValue $ws 'D44' '0.008439'
Set-TextValue $ws 'G44' '20'
Set-TextValue $ws 'D45' '0.00005595'
Set-TextValue $ws 'G45' '20'
Set-TextValue $ws 'G46' '20'
Set-TextValue $ws 'D47' '0.4503'
Set-TextValue $ws 'G47' '20'
Set-TextValue $ws 'D48' '0.1570'
Set-TextValue $ws 'G48' '20'
Set-TextValue $ws 'G49' '20'
Set-TextValue $ws 'D50' '0.01011'
Set-TextValue $ws 'G50' '20'
Set-TextValue $ws 'G51' '20'
